$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.897.80"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.812.83"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.30"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07342"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8675"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.30"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.767.96"
$ws.Range("E12").Value = "  -5.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.380"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07080"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.505"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.68"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008689"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "26.905.84"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "2.035.59"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.42"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.150"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.265"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08930"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.484"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.913"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01945"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.980"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.191"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5290"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.283"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1652"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.389"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4860"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.657"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.95"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06290"
$ws.Range("E51").Value = "  +0.17%  "
